{"js": "// Update .NET MAUI RC2 references to RC3, and the matching VS2022 preview\n// version string (\"VS2022 17.2 Preview 5.0\") to (\"VS2022 17.3 Preview 1.0\").\n\n// 1) Replace every standalone \"RC2\" with \"RC3\" (4 occurrences in the doc).\nconst rc2Results = context.document.body.search(\"RC2\", { matchCase: true });\nrc2Results.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < rc2Results.items.length; i++) {\n  rc2Results.items[i].insertText(\"RC3\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Bump the VS2022 minor version: \"17.2\" -> \"17.3\".\nconst verResults = context.document.body.search(\"17.2\", { matchCase: true });\nverResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < verResults.items.length; i++) {\n  verResults.items[i].insertText(\"17.3\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 3) Bump the Preview build number: \"Preview 5\" -> \"Preview 1\".\nconst previewResults = context.document.body.search(\"Preview 5\", { matchCase: true });\npreviewResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < previewResults.items.length; i++) {\n  previewResults.items[i].insertText(\"Preview 1\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update .NET MAUI RC2 references to RC3, and the matching VS2022 preview\n# version string (\"VS2022 17.2 Preview 5.0\") to (\"VS2022 17.3 Preview 1.0\").\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replaceText) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $findText\n  $find.Replacement.Text = $replaceText\n  $find.Forward = $true\n  $find.Wrap = 1              # wdFindContinue - keep going across the whole story\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.Execute([ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,2) | Out-Null   # 2 = wdReplaceAll\n}\n\n# 1) \"RC2\" -> \"RC3\" (4 occurrences in the document).\nReplace-AllText \"RC2\" \"RC3\"\n\n# 2) VS2022 minor version bump: \"17.2\" -> \"17.3\".\nReplace-AllText \"17.2\" \"17.3\"\n\n# 3) Preview build bump: \"Preview 5\" -> \"Preview 1\".\nReplace-AllText \"Preview 5\" \"Preview 1\"\n"}
